$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "34.050.85"
$ws.Cells.Item(2, 5).Value = "  -1.81%  "
$ws.Cells.Item(3, 4).Value = "1.787.26"
$ws.Cells.Item(3, 5).Value = "  -1.15%  "
$ws.Cells.Item(4, 5).Value = "  +0.12%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "222.95"
$ws.Cells.Item(5, 5).Value = "  -0.99%  "
$ws.Cells.Item(6, 5).Value = "  -0.52%  "
$ws.Cells.Item(7, 5).Value = "  +0.17%  "
$ws.Cells.Item(8, 5).Value = "  -0.86%  "
$ws.Cells.Item(9, 5).Value = "  -1.97%  "
$ws.Cells.Item(10, 5).Value = "  -0.35%  "
$ws.Cells.Item(11, 5).Value = "  +0.24%  "
$ws.Cells.Item(12, 4).Value = "2.045.96"
$ws.Cells.Item(12, 5).Value = "  -1.13%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "11.03"
$ws.Cells.Item(13, 5).Value = "  -0.66%  "
$ws.Cells.Item(14, 4).Value = "1.801.30"
$ws.Cells.Item(14, 5).Value = "  -0.31%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.625"
$ws.Cells.Item(15, 5).Value = "  -3.17%  "
$ws.Cells.Item(16, 4).Value = "34.055.87"
$ws.Cells.Item(16, 5).Value = "  -1.85%  "
$ws.Cells.Item(17, 5).Value = "  -3.93%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "67.92"
$ws.Cells.Item(18, 5).Value = "  -2.54%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "243.91"
$ws.Cells.Item(19, 5).Value = "  -4.15%  "
$ws.Cells.Item(20, 5).Value = "  -2.79%  "
$ws.Cells.Item(21, 5).Value = "  +0.13%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "10.72"
$ws.Cells.Item(22, 5).Value = "  -1.09%  "
$ws.Cells.Item(23, 5).Value = "  -4.30%  "
$ws.Cells.Item(24, 5).Value = "  -2.66%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "158.93"
$ws.Cells.Item(25, 5).Value = "  -1.65%  "
$ws.Cells.Item(26, 5).Value = "  -1.36%  "
$ws.Cells.Item(27, 5).Value = "  -2.18%  "
$ws.Cells.Item(28, 5).Value = "  -2.28%  "
$ws.Cells.Item(29, 5).Value = "  +0.15%  "
$ws.Cells.Item(30, 5).Value = "  -3.00%  "
$ws.Cells.Item(31, 5).Value = "  -0.21%  "
$ws.Cells.Item(32, 5).Value = "  -3.62%  "
$ws.Cells.Item(33, 5).Value = "  -4.59%  "
$ws.Cells.Item(34, 5).Value = "  -5.30%  "
$ws.Cells.Item(35, 4).Value = "1.383.19"
$ws.Cells.Item(35, 5).Value = "  -3.97%  "
$ws.Cells.Item(36, 5).Value = "  +0.53%  "
$ws.Cells.Item(37, 5).Value = "  -1.73%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "79.68"
$ws.Cells.Item(39, 5).Value = "  -6.40%  "
$ws.Cells.Item(40, 5).Value = "  +0.35%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.912"
$ws.Cells.Item(41, 5).Value = "  -5.57%  "
$ws.Cells.Item(42, 5).Value = "  -3.93%  "
$ws.Cells.Item(43, 5).Value = "  +0.50%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.0500"
$ws.Cells.Item(44, 5).Value = "  +1.16%  "
$ws.Cells.Item(45, 2).Value = "WEMIXToken"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "1.05"
$ws.Cells.Item(45, 5).Value = "  -0.55%  "
$ws.Cells.Item(46, 2).Value = "FraxShare"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "5.89"
$ws.Cells.Item(46, 5).Value = "  -3.72%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "107.40"
$ws.Cells.Item(47, 5).Value = "  +1.22%  "
$ws.Cells.Item(48, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(48, 4).Value = "0.0₆0134"
$ws.Cells.Item(48, 5).Value = "  +5.77%  "
$ws.Cells.Item(49, 2).Value = "RocketPoolETH"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(49, 4).Value = "1.944.24"
$ws.Cells.Item(49, 5).Value = "  -0.92%  "
$ws.Cells.Item(50, 5).Value = "  -0.01%  "
$ws.Cells.Item(51, 5).Value = "  -1.95%  "
